# Update Tnfsf13-Sdc2 LR-pair sheet with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => { col => value }
$updates = @{
    2  = @{ G=0.5347833333333334; H=1.60435; I=0.196822066153855; J=0.196822066153855;
            M=1.116695; N=3.350085; O=0.008174214292497491; P=0.008174214292497492;
            Q=0.5971898744166667; R=5.374708869750001; S=0.001608865746233728; T=0.001608865746233729 }
    3  = @{ G=0.5347833333333334; H=1.60435; I=0.196822066153855; J=0.196822066153855;
            O=0.8193429796700005; P=0.8193429796700005;
            Q=59.8593716318889; R=538.7343446870001; S=0.1612647781473055; T=0.1612647781473055 }
    4  = @{ G=0.5347833333333334; H=1.60435; I=0.196822066153855; J=0.196822066153855;
            O=0.172482806037502; P=0.1724828060375021;
            Q=12.60120931391667; R=113.41088382525; S=0.03394842226031577; T=0.03394842226031577 }
    5  = @{ I=0.1891972429821067; J=0.1891972429821067;
            M=1.116695; N=3.350085; O=0.008174214292497491; P=0.008174214292497492;
            Q=0.57405493187; R=5.16649438683; S=0.001546538807685457; T=0.001546538807685457 }
    6  = @{ I=0.1891972429821067; J=0.1891972429821067;
            O=0.8193429796700005; P=0.8193429796700005;
            S=0.1550174328103084; T=0.1550174328103084 }
    7  = @{ I=0.1891972429821067; J=0.1891972429821067;
            O=0.172482806037502; P=0.1724828060375021;
            S=0.03263327136411286; T=0.03263327136411286 }
    8  = @{ I=0.6139806908640383; J=0.6139806908640382;
            M=1.116695; N=3.350085; O=0.008174214292497491; P=0.008174214292497492;
            Q=1.862916383495; R=16.766247451455; S=0.005018809738578305; T=0.005018809738578305 }
    9  = @{ I=0.6139806908640383; J=0.6139806908640382;
            O=0.8193429796700005; P=0.8193429796700005;
            S=0.5030607687123866; T=0.5030607687123865 }
    10 = @{ I=0.6139806908640383; J=0.6139806908640382;
            O=0.172482806037502; P=0.1724828060375021;
            S=0.1059011124130734; T=0.1059011124130734 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
